$wb = $excel.ActiveWorkbook

# Rename the sheet tabs with fresh timestamp-based identifiers.
$wb.Worksheets.Item("GNG_TO-16502911495723593").Name = "GNG_TO-16504778010945537"
$wb.Worksheets.Item("NB_TO-16502911521125338").Name = "NB_TO-16504778031855903"
$wb.Worksheets.Item("RS_TO-16502911521135383").Name = "RS_TO-16504778031865525"
$wb.Worksheets.Item("TOL_TO-16502911521770287").Name = "TOL_TO-16504778032335904"
$wb.Worksheets.Item("vSAT_TO-16502911522652702").Name = "vSAT_TO-16504778032975893"

# GNG sheet (formerly GNG_TO-...) stimulus file names.
$wsGng = $wb.Worksheets.Item("GNG_TO-16504778010945537")
$wsGng.Range("B2").Value = "go_stims-16504778010585546.csv"
$wsGng.Range("B3").Value = "GNG_stims-16504778010775535.csv"
$wsGng.Range("B4").Value = "go_stims-16504778010785558.csv"
$wsGng.Range("B5").Value = "GNG_stims-16504778010935874.csv"

# NB sheet (formerly NB_TO-...) stimulus file names.
$wsNb = $wb.Worksheets.Item("NB_TO-16504778031855903")
$wsNb.Range("B2").Value = "TB-1650477803160554.csv"
$wsNb.Range("B3").Value = "TB-16504778027045796.csv"
$wsNb.Range("B4").Value = "OB-16504778021815524.csv"
$wsNb.Range("B5").Value = "ZB-match_3-1650477801426588.csv"
$wsNb.Range("B6").Value = "TB-16504778027315931.csv"
$wsNb.Range("B7").Value = "ZB-match_6-1650477801603587.csv"
$wsNb.Range("B8").Value = "OB-16504778018815887.csv"
$wsNb.Range("B9").Value = "OB-16504778025285552.csv"
$wsNb.Range("B10").Value = "ZB-match_9-16504778013665519.csv"

# TOL sheet (formerly TOL_TO-...) stimulus file names.
$wsTol = $wb.Worksheets.Item("TOL_TO-16504778032335904")
$wsTol.Range("B2").Value = "MM_stims-1650477803201587.csv"
$wsTol.Range("B3").Value = "ZM_stims-1650477803188557.csv"
$wsTol.Range("B4").Value = "MM_stims-16504778032175908.csv"
$wsTol.Range("B5").Value = "ZM_stims-16504778032025545.csv"
$wsTol.Range("B6").Value = "MM_stims-16504778032335904.csv"
$wsTol.Range("B7").Value = "ZM_stims-16504778032175908.csv"

# vSAT sheet (formerly vSAT_TO-...) stimulus file names.
$wsVsat = $wb.Worksheets.Item("vSAT_TO-16504778032975893")
$wsVsat.Range("B2").Value = "SAT_stims-16504778032365549.csv"
$wsVsat.Range("B3").Value = "vSAT_stims-16504778032815893.csv"
$wsVsat.Range("B4").Value = "vSAT_stims-16504778032655928.csv"
$wsVsat.Range("B5").Value = "SAT_stims-165047780324959.csv"
